$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New French rows of master data, continuing the pattern of id/name/descr/lang_code/is_active/cr_by/cr_dtimes
$rows = @(
    @(10013, "Pré-inscription", "Portail Web pour les pré-inscriptions", "fra"),
    @(10014, "Client dinscription", "Application de bureau pour les inscriptions", "fra"),
    @(10015, "Processeur dinscription", "Demande de post-inscription", "fra"),
    @(10016, "Authentification ID", "Application pour lauthentification du fournisseur de services tiers", "fra"),
    @(10017, "Contrôle didentité", "Portail Web pour la configuration dapplications", "fra"),
    @(10018, "Portail Résident", "Portail Web pour les services de génération de post-ID", "fra")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $true
    $ws.Cells.Item($r, 6).Value = "superadmin"
    $ws.Cells.Item($r, 7).Value = "now()"
}

# Column widths to match new layout (column A was auto-fit to its content,
# column B was given a custom width)
$ws.Columns.Item(1).ColumnWidth = 5
$ws.Columns.Item(2).ColumnWidth = 18.5

# Scroll/selection state to match the saved view
$wb.Windows.Item(1).ScrollRow = 10
$ws.Range("A20:XFD1048576").Select() | Out-Null

# Page setup to portrait / paper size 9 (A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
